# Refresh cached Universalis market-price figures on the Golem Profits sheets.
# (generated from scheduled-runner price snapshot; values only, no formulas)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 112.24
$ws.Range("I33").Value = 82.43478399999999
$ws.Range("J33").Value = 455
$ws.Range("K33").Value = 82.43478399999999
$ws.Range("L33").Value = 455
$ws.Range("M33").Value = 146.565216
$ws.Range("N33").Value = -913

$ws.Range("H34").Value = 3661.25
$ws.Range("I34").Value = 3661.25
$ws.Range("K34").Value = 3661.25
$ws.Range("M34").Value = -3458.25

$ws.Range("H36").Value = 3661.25
$ws.Range("I36").Value = 3661.25
$ws.Range("K36").Value = 3661.25
$ws.Range("M36").Value = -2946.25

$ws.Range("H40").Value = 1883.7097
$ws.Range("J40").Value = 5331.6665
$ws.Range("L40").Value = 5331.6665
$ws.Range("N40").Value = -5681.6665

$ws.Range("H95").Value = 41833
$ws.Range("J95").Value = 41833
$ws.Range("L95").Value = 41833
$ws.Range("N95").Value = -47325

$ws.Range("H101").Value = 512.5
$ws.Range("J101").Value = 525
$ws.Range("L101").Value = 1575
$ws.Range("N101").Value = -4819

$ws.Range("H141").Value = 566
$ws.Range("I141").Value = 566
$ws.Range("K141").Value = 1698
$ws.Range("M141").Value = 3482

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 68000
$ws.Range("J106").Value = 68000
$ws.Range("L106").Value = 68000
$ws.Range("N106").Value = -70524

$ws.Range("H119").Value = 62500
$ws.Range("J119").Value = 62500
$ws.Range("L119").Value = 62500
$ws.Range("N119").Value = -72176

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2099.5
$ws.Range("J20").Value = 2099
$ws.Range("L20").Value = 2099
$ws.Range("N20").Value = -2593

$ws.Range("H26").Value = 24263
$ws.Range("I26").Value = 24263
$ws.Range("K26").Value = 24263
$ws.Range("M26").Value = -23971

$ws.Range("H29").Value = 885.7143
$ws.Range("I29").Value = 875
$ws.Range("J29").Value = 900
$ws.Range("K29").Value = 875
$ws.Range("L29").Value = 900
$ws.Range("M29").Value = -586
$ws.Range("N29").Value = -1478

$ws.Range("H86").Value = 8142.4287
$ws.Range("I86").Value = 7499.75
$ws.Range("J86").Value = 8999.333000000001
$ws.Range("K86").Value = 7499.75
$ws.Range("L86").Value = 8999.333000000001
$ws.Range("M86").Value = -6376.75
$ws.Range("N86").Value = -11245.333

$ws.Range("H89").Value = 8142.4287
$ws.Range("I89").Value = 7499.75
$ws.Range("J89").Value = 8999.333000000001
$ws.Range("K89").Value = 37498.75
$ws.Range("L89").Value = 44996.665
$ws.Range("M89").Value = -31882.75
$ws.Range("N89").Value = -56228.665

$ws.Range("H99").Value = 2458.6
$ws.Range("I99").Value = 2685.4119
$ws.Range("K99").Value = 2685.4119
$ws.Range("M99").Value = -1187.4119

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 12959.333
$ws.Range("J15").Value = 12959.333
$ws.Range("L15").Value = 12959.333
$ws.Range("N15").Value = -13299.333

$ws.Range("H50").Value = 43999
$ws.Range("J50").Value = 43999
$ws.Range("L50").Value = 43999
$ws.Range("N50").Value = -45249

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H62").Value = 2466.6667
$ws.Range("I62").Value = 2466.6667
$ws.Range("K62").Value = 2466.6667
$ws.Range("M62").Value = -1842.6667

$ws.Range("H65").Value = 2466.6667
$ws.Range("I65").Value = 2466.6667
$ws.Range("K65").Value = 12333.3335
$ws.Range("M65").Value = -9213.333500000001

$ws.Range("H68").Value = 45147.5
$ws.Range("J68").Value = 45147.5
$ws.Range("L68").Value = 45147.5
$ws.Range("N68").Value = -46645.5

$ws.Range("H71").Value = 45147.5
$ws.Range("J71").Value = 45147.5
$ws.Range("L71").Value = 135442.5
$ws.Range("N71").Value = -142930.5

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2563.1462
$ws.Range("I4").Value = 1896.0358
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 5688.107400000001
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = -5576.107400000001
$ws.Range("N4").Value = -12224

$ws.Range("H12").Value = 187.66667
$ws.Range("J12").Value = 272.83334
$ws.Range("L12").Value = 818.5000200000001
$ws.Range("N12").Value = -1164.50002

$ws.Range("H75").Value = 2360.25
$ws.Range("I75").Value = 95
$ws.Range("J75").Value = 2683.8572
$ws.Range("K75").Value = 285
$ws.Range("L75").Value = 8051.571599999999
$ws.Range("M75").Value = 713
$ws.Range("N75").Value = -10047.5716

$ws.Range("H78").Value = 2360.25
$ws.Range("I78").Value = 95
$ws.Range("J78").Value = 2683.8572
$ws.Range("K78").Value = 855
$ws.Range("L78").Value = 24154.7148
$ws.Range("M78").Value = 4137
$ws.Range("N78").Value = -34138.7148

$ws.Range("H86").Value = 2420.8333
$ws.Range("I86").Value = 1375
$ws.Range("J86").Value = 2943.75
$ws.Range("K86").Value = 4125
$ws.Range("L86").Value = 8831.25
$ws.Range("M86").Value = -2939
$ws.Range("N86").Value = -11203.25

$ws.Range("H89").Value = 2420.8333
$ws.Range("I89").Value = 1375
$ws.Range("J89").Value = 2943.75
$ws.Range("K89").Value = 12375
$ws.Range("L89").Value = 26493.75
$ws.Range("M89").Value = -6447
$ws.Range("N89").Value = -38349.75

$ws.Range("H123").Value = 900
$ws.Range("I123").Value = 900
$ws.Range("K123").Value = 2700
$ws.Range("M123").Value = -250

$ws.Range("H137").Value = 5166.3335
$ws.Range("I137").Value = 2000
$ws.Range("K137").Value = 6000
$ws.Range("M137").Value = -900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 206
$ws.Range("I21").Value = 206
$ws.Range("K21").Value = 206
$ws.Range("M21").Value = -32

$ws.Range("H22").Value = 2537.389
$ws.Range("I22").Value = 1277.7778
$ws.Range("J22").Value = 3797
$ws.Range("K22").Value = 1277.7778
$ws.Range("L22").Value = 3797
$ws.Range("M22").Value = -982.7778000000001
$ws.Range("N22").Value = -4387

$ws.Range("H27").Value = 2537.389
$ws.Range("I27").Value = 1277.7778
$ws.Range("J27").Value = 3797
$ws.Range("K27").Value = 1277.7778
$ws.Range("L27").Value = 3797
$ws.Range("M27").Value = -1170.7778
$ws.Range("N27").Value = -4011

$ws.Range("H46").Value = 226621.33
$ws.Range("I46").Value = 667200
$ws.Range("K46").Value = 667200
$ws.Range("M46").Value = -667012

$ws.Range("H82").Value = 1803.3
$ws.Range("I82").Value = 1761.8572
$ws.Range("J82").Value = 1900
$ws.Range("K82").Value = 1761.8572
$ws.Range("L82").Value = 1900
$ws.Range("M82").Value = -1400.8572
$ws.Range("N82").Value = -2622

$ws.Range("H85").Value = 1803.3
$ws.Range("I85").Value = 1761.8572
$ws.Range("J85").Value = 1900
$ws.Range("K85").Value = 1761.8572
$ws.Range("L85").Value = 1900
$ws.Range("M85").Value = -513.8571999999999
$ws.Range("N85").Value = -4396

$ws.Range("H93").Value = 1449.3
$ws.Range("I93").Value = 1742.5714
$ws.Range("J93").Value = 765
$ws.Range("K93").Value = 1742.5714
$ws.Range("L93").Value = 765
$ws.Range("M93").Value = -494.5714
$ws.Range("N93").Value = -3261

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
